# Insert a new data row at row 14 (pushing existing rows 14-72 down to 15-73)
# and populate it with a new price-report record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 14; formatting of the row below
# (the old row 14) carries over to the new blank row, matching the
# date-formatted style already used for column D.
$ws.Rows.Item(14).Insert()

# Populate the newly inserted row 14 with the new record's data.
$ws.Range("A14").Value = 2
$ws.Range("B14").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C14").Value = "Coquimbo"
$ws.Range("D14").Value = 44608
$ws.Range("E14").Value = 4
$ws.Range("F14").Value = 100112030
$ws.Range("G14").Value = "Poroto granado"
$ws.Range("H14").Value = "Sin especificar"
$ws.Range("I14").Value = "Primera"
$ws.Range("J14").Value = 600
$ws.Range("K14").Value = 22000
$ws.Range("L14").Value = 24000
$ws.Range("M14").Value = 23000
$ws.Range("N14").Value = "$/malla 25 kilos"
$ws.Range("O14").Value = "Provincia de Limarí"
$ws.Range("P14").Value = 920
$ws.Range("Q14").Value = 25
$ws.Range("R14").Value = "Hortaliza"
